$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 55, pushing existing rows 55..106 down to 56..107
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new weekly record
$ws.Cells.Item(55, 1).Value = 10
$ws.Cells.Item(55, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(55, 3).Value = "La Araucanía"
$ws.Cells.Item(55, 4).Value = 44893
$ws.Cells.Item(55, 5).Value = 9
$ws.Cells.Item(55, 6).Value = "Fruta"
$ws.Cells.Item(55, 7).Value = 100101
$ws.Cells.Item(55, 8).Value = "Berries"
$ws.Cells.Item(55, 9).Value = 100101001
$ws.Cells.Item(55, 10).Value = "Arándano (blue)"
$ws.Cells.Item(55, 11).Value = "Sin especificar"
$ws.Cells.Item(55, 12).Value = "Primera"
$ws.Cells.Item(55, 13).Value = 680
$ws.Cells.Item(55, 14).Value = 2600
$ws.Cells.Item(55, 15).Value = 2700
$ws.Cells.Item(55, 16).Value = 2644
$ws.Cells.Item(55, 17).Value = "$/kilo"
$ws.Cells.Item(55, 18).Value = "Región del Maule"
$ws.Cells.Item(55, 19).Value = 2644
$ws.Cells.Item(55, 20).Value = 1
